$wb = $excel.ActiveWorkbook

# Add the new "edges" worksheet as the last sheet in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "edges"

# Column widths to match the authored sheet (~14.43 / ~14.29 chars).
$ws.Columns.Item(1).ColumnWidth = 13.666666666666666
$ws.Columns.Item(2).ColumnWidth = 13.5

# Header row.
$ws.Range("A1").Value = "nodo_a"
$ws.Range("B1").Value = "nodo_b"

# Edge list rows.
$ws.Range("A2").Value = "rain"
$ws.Range("B2").Value = "train"

$ws.Range("A3").Value = "rain"
$ws.Range("B3").Value = "maintenance"

$ws.Range("A4").Value = "maintenance"
$ws.Range("B4").Value = "train"

$ws.Range("A5").Value = "train"
$ws.Range("B5").Value = "appointment"

# Make the new sheet the active/selected tab, like the authored file.
$ws.Activate()
$ws.Range("B7").Select()
